$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 58

# Column A holds date-like text ("2026/01/07"). Pre-format as Text so the
# value is stored as a literal string instead of being auto-parsed into a
# date serial number, then clear the format back down and re-apply the
# same centered alignment used by the rest of the table (rows 3-57) so the
# new row lands on the same cell style.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2026/01/07"
$ws.Cells.Item($row, 2).Value = "逃离鸭科夫"
$ws.Cells.Item($row, 3).Value = 1140

$rng = $ws.Range("A$($row):C$($row)")
$rng.ClearFormats()
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108
